{"js": "// Replace each \"NNxNN=\" multiplication prompt in the document's table\n// cells with its new value, per the commit's regenerated problem set.\n// Each old value occurs exactly once in the document, so a scoped\n// case-sensitive search-and-replace per pair is unambiguous and leaves\n// all run formatting (font/size) untouched.\nconst replacements = [\n  [\"14\u00d761=\", \"33\u00d718=\"],\n  [\"13\u00d748=\", \"12\u00d796=\"],\n  [\"18\u00d774=\", \"17\u00d738=\"],\n  [\"70\u00d731=\", \"74\u00d737=\"],\n  [\"48\u00d743=\", \"80\u00d785=\"],\n  [\"40\u00d780=\", \"99\u00d797=\"],\n  [\"29\u00d747=\", \"16\u00d780=\"],\n  [\"85\u00d792=\", \"85\u00d770=\"],\n  [\"79\u00d780=\", \"65\u00d748=\"],\n  [\"59\u00d764=\", \"78\u00d762=\"],\n  [\"29\u00d754=\", \"38\u00d792=\"],\n  [\"15\u00d742=\", \"42\u00d797=\"],\n  [\"89\u00d733=\", \"62\u00d726=\"],\n  [\"20\u00d759=\", \"18\u00d785=\"],\n  [\"84\u00d786=\", \"69\u00d725=\"],\n  [\"55\u00d734=\", \"87\u00d773=\"],\n  [\"72\u00d757=\", \"37\u00d793=\"],\n  [\"42\u00d799=\", \"84\u00d756=\"],\n  [\"94\u00d731=\", \"72\u00d730=\"],\n  [\"73\u00d715=\", \"74\u00d737=\"],\n  [\"52\u00d743=\", \"77\u00d785=\"],\n  [\"30\u00d746=\", \"44\u00d784=\"],\n  [\"77\u00d756=\", \"12\u00d730=\"],\n  [\"48\u00d764=\", \"15\u00d798=\"],\n  [\"71\u00d719=\", \"75\u00d718=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"NNxNN=\" multiplication prompt in the document's table\n# cells with its new value, per the commit's regenerated problem set.\n# Each old value occurs exactly once in the document, so a whole-document\n# Find/Replace per pair is unambiguous and leaves run formatting intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"14\u00d761=\", \"33\u00d718=\"),\n    @(\"13\u00d748=\", \"12\u00d796=\"),\n    @(\"18\u00d774=\", \"17\u00d738=\"),\n    @(\"70\u00d731=\", \"74\u00d737=\"),\n    @(\"48\u00d743=\", \"80\u00d785=\"),\n    @(\"40\u00d780=\", \"99\u00d797=\"),\n    @(\"29\u00d747=\", \"16\u00d780=\"),\n    @(\"85\u00d792=\", \"85\u00d770=\"),\n    @(\"79\u00d780=\", \"65\u00d748=\"),\n    @(\"59\u00d764=\", \"78\u00d762=\"),\n    @(\"29\u00d754=\", \"38\u00d792=\"),\n    @(\"15\u00d742=\", \"42\u00d797=\"),\n    @(\"89\u00d733=\", \"62\u00d726=\"),\n    @(\"20\u00d759=\", \"18\u00d785=\"),\n    @(\"84\u00d786=\", \"69\u00d725=\"),\n    @(\"55\u00d734=\", \"87\u00d773=\"),\n    @(\"72\u00d757=\", \"37\u00d793=\"),\n    @(\"42\u00d799=\", \"84\u00d756=\"),\n    @(\"94\u00d731=\", \"72\u00d730=\"),\n    @(\"73\u00d715=\", \"74\u00d737=\"),\n    @(\"52\u00d743=\", \"77\u00d785=\"),\n    @(\"30\u00d746=\", \"44\u00d784=\"),\n    @(\"77\u00d756=\", \"12\u00d730=\"),\n    @(\"48\u00d764=\", \"15\u00d798=\"),\n    @(\"71\u00d719=\", \"75\u00d718=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
